$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp (A1) -------------------------------
$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 17:41"

# --- Swap Grecia / Noruega ranking (rows 91 & 92) and refresh their stats -
# Row 91 was Noruega, now becomes Grecia (Grecia overtook Noruega in total cases)
$ws.Range("A91").Value = "Grecia"
$ws.Range("B91").Value = 11200
$ws.Range("C91").Value = 202
$ws.Range("D91").Value = 3804
$ws.Range("E91").Value = 7117
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 279

# Row 92 was Grecia, now becomes Noruega
$ws.Range("A92").Value = "Noruega"
$ws.Range("B92").Value = 11160
$ws.Range("C92").Value = 40
$ws.Range("D92").Value = 9348
$ws.Range("E92").Value = 1548
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 264

# --- Refresh daily COVID counters for the other updated countries --------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 6339403
$ws.Range("C4").Value = 4159
$ws.Range("E4").Value = 2572345
$ws.Range("G4").Value = 134
$ws.Range("H4").Value = 191192

# Row 52: Singapur
$ws.Range("D52").Value = 56174
$ws.Range("E52").Value = 747

# Row 65: Moldavia
$ws.Range("B65").Value = 38906
$ws.Range("C65").Value = 534
$ws.Range("E65").Value = 10842
$ws.Range("G65").Value = 11
$ws.Range("H65").Value = 1047

# Row 68: Kenia
$ws.Range("B68").Value = 34884
$ws.Range("C68").Value = 179
$ws.Range("D68").Value = 21059
$ws.Range("E68").Value = 13236
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 589

# Row 94: Albania
$ws.Range("B94").Value = 9967
$ws.Range("C94").Value = 123
$ws.Range("D94").Value = 5882
$ws.Range("E94").Value = 3779
$ws.Range("G94").Value = 5
$ws.Range("H94").Value = 306

# Row 127: Somalia
$ws.Range("B127").Value = 3332
$ws.Range("C127").Value = 22
$ws.Range("D127").Value = 2696
$ws.Range("E127").Value = 539

# Row 129: Sri Lanka
$ws.Range("B129").Value = 3115
$ws.Range("C129").Value = 4
$ws.Range("E129").Value = 196

# Row 176: Papua Nueva Guinea
$ws.Range("B176").Value = 479
$ws.Range("C176").Value = 8
$ws.Range("E176").Value = 242
